# Updates cryptos list: refreshed Price (D) / Volume(1h) (E) figures for Thu Oct 31
# 2024 GitHub Actions run; rows 37/38 and 42/43 also swap rank order (Coin/Link/Price/Volume).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.368.42"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "2.640.73"
$ws.Range("E3").Value = "  -1.55%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'583.12"
$ws.Range("E5").Value = "  -3.22%  "

$ws.Range("D6").Value = "'175.11"
$ws.Range("E6").Value = "  -1.83%  "

$ws.Range("D8").Value = "'0.520"
$ws.Range("E8").Value = "  -1.09%  "

$ws.Range("D9").Value = "2.639.72"
$ws.Range("E9").Value = "  -1.54%  "

$ws.Range("E10").Value = "  -0.09%  "

$ws.Range("E11").Value = "  +0.93%  "

$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("E13").Value = "  -2.88%  "

$ws.Range("D14").Value = "3.125.29"
$ws.Range("E14").Value = "  -1.43%  "

$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").Value = "72.217.99"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").Value = "'25.84"
$ws.Range("E17").Value = "  -2.14%  "

$ws.Range("D18").Value = "2.653.37"
$ws.Range("E18").Value = "  -1.38%  "

$ws.Range("D19").Value = "'8.43"
$ws.Range("E19").Value = "  +4.80%  "

$ws.Range("D20").Value = "'12.10"
$ws.Range("E20").Value = "  +1.31%  "

$ws.Range("D21").Value = "'373.65"
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").Value = "'4.11"
$ws.Range("E22").Value = "  -1.72%  "

$ws.Range("D23").Value = "'2.03"
$ws.Range("E23").Value = "  -0.88%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").Value = "'70.89"
$ws.Range("E25").Value = "  -2.16%  "

$ws.Range("E26").Value = "  -2.63%  "

$ws.Range("D27").Value = "'9.46"
$ws.Range("E27").Value = "  -4.70%  "

$ws.Range("D28").Value = "2.777.53"
$ws.Range("E28").Value = "  -1.48%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "0.0₃0948"
$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("E31").Value = "  -1.78%  "

$ws.Range("D32").Value = "'494.74"
$ws.Range("E32").Value = "  -4.66%  "

$ws.Range("E33").Value = "  -3.14%  "

$ws.Range("E34").Value = "  -2.01%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").Value = "'162.72"
$ws.Range("E36").Value = "  -1.16%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'19.17"
$ws.Range("E37").Value = "  -2.02%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.114"
$ws.Range("E38").Value = "  +4.32%  "

$ws.Range("D39").Value = "'18.86"
$ws.Range("E39").Value = "  -1.45%  "

$ws.Range("E40").Value = "  -2.13%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.71"
$ws.Range("E42").Value = "  -6.67%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.56"
$ws.Range("E43").Value = "  -1.33%  "

$ws.Range("D44").Value = "'4.90"
$ws.Range("E44").Value = "  -3.12%  "

$ws.Range("E45").Value = "  -2.50%  "

$ws.Range("D46").Value = "'39.03"
$ws.Range("E46").Value = "  -0.52%  "

$ws.Range("D47").Value = "'151.83"
$ws.Range("E47").Value = "  -1.61%  "

$ws.Range("D48").Value = "'3.64"
$ws.Range("E48").Value = "  -2.58%  "

$ws.Range("D49").Value = "'0.544"
$ws.Range("E49").Value = "  -1.37%  "

$ws.Range("E50").Value = "  -3.79%  "

$ws.Range("E51").Value = "  -0.88%  "
